$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 3744.9285
$ws.Range("J6").Value = 233.33333
$ws.Range("L6").Value = 699.99999
$ws.Range("N6").Value = -923.99999
$ws.Range("H8").Value = 562.4286
$ws.Range("I8").Value = 572.8333
$ws.Range("J8").Value = 500
$ws.Range("K8").Value = 1718.4999
$ws.Range("L8").Value = 1500
$ws.Range("M8").Value = -1579.4999
$ws.Range("N8").Value = -1778
$ws.Range("H38").Value = 1200.3334
$ws.Range("I38").Value = 123.8
$ws.Range("J38").Value = 1668.3914
$ws.Range("K38").Value = 371.4
$ws.Range("L38").Value = 5005.174199999999
$ws.Range("M38").Value = 0.6000000000000227
$ws.Range("N38").Value = -5749.174199999999
$ws.Range("H39").Value = 555.4583
$ws.Range("I39").Value = 567.3570999999999
$ws.Range("J39").Value = 538.8
$ws.Range("K39").Value = 1702.0713
$ws.Range("L39").Value = 1616.4
$ws.Range("M39").Value = -1406.0713
$ws.Range("N39").Value = -2208.4
$ws.Range("H43").Value = 4274092
$ws.Range("I43").Value = 454.1
$ws.Range("K43").Value = 454.1
$ws.Range("M43").Value = -385.1
$ws.Range("H70").Value = 2720.7058
$ws.Range("I70").Value = 2720.7058
$ws.Range("K70").Value = 8162.117400000001
$ws.Range("M70").Value = -7892.117400000001
$ws.Range("H73").Value = 2720.7058
$ws.Range("I73").Value = 2720.7058
$ws.Range("K73").Value = 8162.117400000001
$ws.Range("M73").Value = -7226.117400000001
$ws.Range("H94").Value = 4333.3335
$ws.Range("I94").Value = 4333.3335
$ws.Range("K94").Value = 4333.3335
$ws.Range("M94").Value = -3882.3335
$ws.Range("H127").Value = 1041.875
$ws.Range("I127").Value = 529.2
$ws.Range("J127").Value = 1896.3334
$ws.Range("K127").Value = 1587.6
$ws.Range("L127").Value = 5689.0002
$ws.Range("M127").Value = 3372.4
$ws.Range("N127").Value = -15609.0002
$ws.Range("H129").Value = 885.05457
$ws.Range("J129").Value = 899.6415
$ws.Range("L129").Value = 2698.9245
$ws.Range("N129").Value = -12698.9245
$ws.Range("H132").Value = 16674365
$ws.Range("I132").Value = 23814314
$ws.Range("J132").Value = 14483.333
$ws.Range("K132").Value = 71442942
$ws.Range("L132").Value = 43449.999
$ws.Range("M132").Value = -71440412
$ws.Range("N132").Value = -48509.999
$ws.Range("H138").Value = 2275.32
$ws.Range("J138").Value = 2275.32
$ws.Range("L138").Value = 6825.960000000001
$ws.Range("N138").Value = -17105.96

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 169.5
$ws.Range("I5").Value = 169.5
$ws.Range("K5").Value = 169.5
$ws.Range("M5").Value = -57.5
$ws.Range("H32").Value = 2933.1836
$ws.Range("I32").Value = 3079.7173
$ws.Range("K32").Value = 3079.7173
$ws.Range("M32").Value = -2792.7173
$ws.Range("H45").Value = 2097.0833
$ws.Range("I45").Value = 2116.5
$ws.Range("K45").Value = 2116.5
$ws.Range("M45").Value = -1739.5
$ws.Range("H61").Value = 1772.6428
$ws.Range("I61").Value = 1445.7273
$ws.Range("K61").Value = 1445.7273
$ws.Range("M61").Value = -1233.7273
$ws.Range("H132").Value = 2993.5652
$ws.Range("I132").Value = 2769.6667
$ws.Range("J132").Value = 3799.6
$ws.Range("K132").Value = 8309.000100000001
$ws.Range("L132").Value = 11398.8
$ws.Range("M132").Value = -5779.000100000001
$ws.Range("N132").Value = -16458.8
$ws.Range("H136").Value = 1772.6428
$ws.Range("I136").Value = 1445.7273
$ws.Range("K136").Value = 4337.1819
$ws.Range("M136").Value = -1787.1819

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 169.5
$ws.Range("I4").Value = 169.5
$ws.Range("K4").Value = 169.5
$ws.Range("M4").Value = -54.5
$ws.Range("H86").Value = 3374.7144
$ws.Range("I86").Value = 3561.524
$ws.Range("J86").Value = 2814.2856
$ws.Range("K86").Value = 3561.524
$ws.Range("L86").Value = 2814.2856
$ws.Range("M86").Value = -2438.524
$ws.Range("N86").Value = -5060.2856
$ws.Range("H89").Value = 3374.7144
$ws.Range("I89").Value = 3561.524
$ws.Range("J89").Value = 2814.2856
$ws.Range("K89").Value = 17807.62
$ws.Range("L89").Value = 14071.428
$ws.Range("M89").Value = -12191.62
$ws.Range("N89").Value = -25303.428
$ws.Range("H94").Value = 22728252
$ws.Range("I94").Value = 31250620
$ws.Range("K94").Value = 31250620
$ws.Range("M94").Value = -31250169
$ws.Range("H105").Value = 111114070
$ws.Range("I105").Value = 142859950
$ws.Range("J105").Value = 3500
$ws.Range("K105").Value = 142859950
$ws.Range("L105").Value = 3500
$ws.Range("M105").Value = -142858203
$ws.Range("N105").Value = -6994
$ws.Range("H107").Value = 1515.9333
$ws.Range("I107").Value = 1192.7
$ws.Range("K107").Value = 1192.7
$ws.Range("M107").Value = 727.3
$ws.Range("H134").Value = 7286.263
$ws.Range("I134").Value = 1989.3334
$ws.Range("K134").Value = 5968.0002
$ws.Range("M134").Value = -3433.0002

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 425.9
$ws.Range("J7").Value = 462.5
$ws.Range("L7").Value = 462.5
$ws.Range("N7").Value = -688.5
$ws.Range("H86").Value = 4778699
$ws.Range("I86").Value = 8335770
$ws.Range("K86").Value = 8335770
$ws.Range("M86").Value = -8334647
$ws.Range("H89").Value = 4778699
$ws.Range("I89").Value = 8335770
$ws.Range("K89").Value = 41678850
$ws.Range("M89").Value = -41673234
$ws.Range("H94").Value = 688.6923
$ws.Range("I94").Value = 1200
$ws.Range("K94").Value = 1200
$ws.Range("M94").Value = -749
$ws.Range("H99").Value = 2610
$ws.Range("I99").Value = 2512
$ws.Range("J99").Value = 2757
$ws.Range("K99").Value = 2512
$ws.Range("L99").Value = 2757
$ws.Range("M99").Value = -1014
$ws.Range("N99").Value = -5753
$ws.Range("H103").Value = 1683.1666
$ws.Range("I103").Value = 1683.1666
$ws.Range("K103").Value = 1683.1666
$ws.Range("M103").Value = -511.1666
$ws.Range("H126").Value = 2610
$ws.Range("I126").Value = 2512
$ws.Range("J126").Value = 2757
$ws.Range("K126").Value = 7536
$ws.Range("L126").Value = 8271
$ws.Range("M126").Value = -5066
$ws.Range("N126").Value = -13211

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 25001408
$ws.Range("J131").Value = 1584.875
$ws.Range("L131").Value = 4754.625
$ws.Range("N131").Value = -14834.625
$ws.Range("H140").Value = 30537.5
$ws.Range("I140").Value = 34709.97
$ws.Range("J140").Value = 2999.2
$ws.Range("K140").Value = 104129.91
$ws.Range("L140").Value = 8997.599999999999
$ws.Range("M140").Value = -98949.91
$ws.Range("N140").Value = -19357.6

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H5").Value = 4000
$ws.Range("J5").Value = 4000
$ws.Range("L5").Value = 4000
$ws.Range("N5").Value = -4226
$ws.Range("H132").Value = 79946.234
$ws.Range("I132").Value = 2949.6667
$ws.Range("K132").Value = 8849.000100000001
$ws.Range("M132").Value = -6319.000100000001
$ws.Range("H136").Value = 2077.6
$ws.Range("I136").Value = 2011.3572
$ws.Range("J136").Value = 3005
$ws.Range("K136").Value = 6034.071599999999
$ws.Range("L136").Value = 9015
$ws.Range("M136").Value = -3484.071599999999
$ws.Range("N136").Value = -14115

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H104").Value = 0
$ws.Range("J104").Value = 0
$ws.Range("L104").Value = 0
$ws.Range("N104").ClearContents()
